$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: update label text (values only changed, no structural change)
$ws.Range("A1").Value = "legenda"
$ws.Range("B1").Value = "area"
$ws.Range("D1").Value = "area_km2"

# Data rows: B2/B3 numeric updates
$ws.Range("B2").Value = 74207.74782040001
$ws.Range("B3").Value = 100191.517033
